$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the contents of a whole paragraph (everything except its
# trailing paragraph mark) with an explicit run-level OOXML fragment via
# Range.InsertXML. Supplying the *entire* paragraph content (rather than
# just the edited tail) avoids the host's tendency to coalesce the
# untouched leading runs (e.g. collapsing a <w:tab/> run into a literal
# tab character) when only part of a paragraph is replaced.
# ---------------------------------------------------------------------------
function Set-ParagraphXml {
    param($Paragraph, $InnerXml)
    $startPos = $Paragraph.Range.Start
    $endPos = $Paragraph.Range.End - 1
    $target = $d.Range($startPos, $endPos)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $InnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) Paragraph 1 - wrap "SourceTree" in spell-check proofErr markers and
#    split it (and the following period) into their own runs.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Xml = '<w:p>' +
    '<w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">20/07/2017 16:07 </w:t></w:r>' +
    '<w:r><w:tab/><w:t xml:space="preserve">Files have been downloaded and added to a GitHub repository. Hosted on localhost through Node.js and managed using </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>SourceTree</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Sublime Text used for IDE.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p1 $p1Xml

# ---------------------------------------------------------------------------
# 2) Paragraph 6 - wrap "OnClick" in spell-check proofErr markers and split
#    the sentence around it into separate runs.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$p6Xml = '<w:p>' +
    '<w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' +
    '<w:r><w:t>20/07/2017 19:17</w:t></w:r>' +
    '<w:r><w:tab/><w:t xml:space="preserve">Managed to create Handlebars template for accordion element with content for each section. Came across problem where JavaScript was unable to loop through array of sections and add </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>OnClick</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> event allowing accordion to function properly.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p6 $p6Xml

# ---------------------------------------------------------------------------
# 3) Final paragraph (currently paragraph 8) - drop the trailing
#    bookmarkStart/bookmarkEnd pair (they move to the very end of the
#    document) and append two brand-new log entries: a "Break" entry and
#    the "Completed Advanced task 'CSS boxes'." entry, which now carries
#    the bookmark pair.
# ---------------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastXml = '<w:p>' +
    '<w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' +
    '<w:r><w:t>20/07/2017 20:05</w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t>Completed Advanced task &#8216;</w:t></w:r>' +
    '<w:r><w:t>Convert the CSS to LESS or SASS</w:t></w:r>' +
    '<w:r><w:t>&#8217;.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Converted CSS to LESS.</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
    '<w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' +
    '<w:r><w:t>20/07/2017 20:40</w:t></w:r>' +
    '<w:r><w:tab/><w:t>Break</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
    '<w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' +
    '<w:r><w:t>20/07/2017 22:</w:t></w:r>' +
    '<w:r><w:t>00</w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t>Completed Advanced task &#8216;</w:t></w:r>' +
    '<w:r><w:t>CSS boxes</w:t></w:r>' +
    '<w:r><w:t>&#8217;.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
Set-ParagraphXml $pLast $lastXml

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
